$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("10:15").Insert()
